$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: GNK Dinamo Zagreb
$ws.Cells.Item(2,1).Value = "GNK Dinamo Zagreb ✓ - NK Osijek: 2:1"
$ws.Cells.Item(2,2).Value = "GNK Dinamo Zagreb"
$ws.Cells.Item(2,3).Value = 72
$ws.Cells.Item(2,4).Value = 75
$ws.Cells.Item(2,5).Value = $null
$ws.Cells.Item(2,6).Value = 1.36
$ws.Cells.Item(2,7).Value = "✓"

# Row 3: Bayern Munich
$ws.Cells.Item(3,1).Value = "Bayern Munich ✓ - Borussia Dortmund: 2:1"
$ws.Cells.Item(3,2).Value = "Bayern Munich"
$ws.Cells.Item(3,3).Value = 72
$ws.Cells.Item(3,4).Value = 73
$ws.Cells.Item(3,5).Value = 83
$ws.Cells.Item(3,6).Value = 1.4
$ws.Cells.Item(3,7).Value = "✓"

# Row 4: PSV Eindhoven
$ws.Cells.Item(4,1).Value = "PSV Eindhoven ✓ - Go Ahead Eagles: 2:1"
$ws.Cells.Item(4,2).Value = "PSV Eindhoven"
$ws.Cells.Item(4,3).Value = 71
$ws.Cells.Item(4,4).Value = 91
$ws.Cells.Item(4,5).Value = $null
$ws.Cells.Item(4,6).Value = 1.25
$ws.Cells.Item(4,7).Value = "✓"

# Row 5: FC Basel 1893
$ws.Cells.Item(5,1).Value = "FC Basel 1893 ✓ - FC Winterthur: 3:0"
$ws.Cells.Item(5,2).Value = "FC Basel 1893"
$ws.Cells.Item(5,3).Value = 71
$ws.Cells.Item(5,4).Value = 88
$ws.Cells.Item(5,5).Value = $null
$ws.Cells.Item(5,6).Value = 1.27
$ws.Cells.Item(5,7).Value = "✓"

# Row 6: AS Saint-Étienne
$ws.Cells.Item(6,1).Value = "AS Saint-Étienne X - Le Mans FC: 2:3"
$ws.Cells.Item(6,2).Value = "AS Saint-Étienne"
$ws.Cells.Item(6,3).Value = 71
$ws.Cells.Item(6,4).Value = 100
$ws.Cells.Item(6,5).Value = $null
$ws.Cells.Item(6,6).Value = 1.45
$ws.Cells.Item(6,7).Value = "X"

# Row 7: Sporting Club Escaldes
$ws.Cells.Item(7,1).Value = "UE Santa Coloma B - Sporting Club Escaldes ✓: 0:3"
$ws.Cells.Item(7,2).Value = "Sporting Club Escaldes"
$ws.Cells.Item(7,3).Value = 70
$ws.Cells.Item(7,4).Value = 86
$ws.Cells.Item(7,5).Value = $null
$ws.Cells.Item(7,6).Value = 1.67
$ws.Cells.Item(7,7).Value = "✓"

# Row 8: Västerås SK
$ws.Cells.Item(8,1).Value = "Västerås SK ✓ - GIF Sundsvall: 2:1"
$ws.Cells.Item(8,2).Value = "Västerås SK"
$ws.Cells.Item(8,3).Value = 70
$ws.Cells.Item(8,4).Value = 85
$ws.Cells.Item(8,5).Value = $null
$ws.Cells.Item(8,6).Value = 1.62
$ws.Cells.Item(8,7).Value = "✓"

# Row 9: Galatasaray
$ws.Cells.Item(9,1).Value = "Basaksehir FK - Galatasaray ✓: 1:2"
$ws.Cells.Item(9,2).Value = "Galatasaray"
$ws.Cells.Item(9,3).Value = 70
$ws.Cells.Item(9,4).Value = 82
$ws.Cells.Item(9,5).Value = $null
$ws.Cells.Item(9,6).Value = 1.7
$ws.Cells.Item(9,7).Value = "✓"

# Row 10: Thep Xanh Nam Dinh FC
$ws.Cells.Item(10,1).Value = "Thep Xanh Nam Dinh FC X - Becamex Ho Chi Minh City FC: 1:2"
$ws.Cells.Item(10,2).Value = "Thep Xanh Nam Dinh FC"
$ws.Cells.Item(10,3).Value = 67
$ws.Cells.Item(10,4).Value = 75
$ws.Cells.Item(10,5).Value = $null
$ws.Cells.Item(10,6).Value = 1.36
$ws.Cells.Item(10,7).Value = "X"

# Row 11: SK Slavia Prague
$ws.Cells.Item(11,1).Value = "SK Slavia Prague  - FC Zlin: 0:0"
$ws.Cells.Item(11,2).Value = "SK Slavia Prague"
$ws.Cells.Item(11,3).Value = 66
$ws.Cells.Item(11,4).Value = 75
$ws.Cells.Item(11,5).Value = $null
$ws.Cells.Item(11,6).Value = 1.18
$ws.Cells.Item(11,7).Value = $null

# Row 12: FC Barcelona
$ws.Cells.Item(12,1).Value = "FC Barcelona ✓ - Girona FC: 2:1"
$ws.Cells.Item(12,2).Value = "FC Barcelona"
$ws.Cells.Item(12,3).Value = 65
$ws.Cells.Item(12,4).Value = $null
$ws.Cells.Item(12,5).Value = 83
$ws.Cells.Item(12,6).Value = 1.18
$ws.Cells.Item(12,7).Value = "✓"

# Row 13: Arsenal FC
$ws.Cells.Item(13,1).Value = "Fulham FC - Arsenal FC ✓: 0:1"
$ws.Cells.Item(13,2).Value = "Arsenal FC"
$ws.Cells.Item(13,3).Value = 60
$ws.Cells.Item(13,4).Value = 78
$ws.Cells.Item(13,5).Value = 76
$ws.Cells.Item(13,6).Value = 1.52
$ws.Cells.Item(13,7).Value = "✓"

# Row 14: Al-Ain FC
$ws.Cells.Item(14,1).Value = "Al-Ain FC ✓ - FC Baniyas: 4:0"
$ws.Cells.Item(14,2).Value = "Al-Ain FC"
$ws.Cells.Item(14,3).Value = 59
$ws.Cells.Item(14,4).Value = 100
$ws.Cells.Item(14,5).Value = $null
$ws.Cells.Item(14,6).Value = 1.45
$ws.Cells.Item(14,7).Value = "✓"

# Row 15: Club Nacional
$ws.Cells.Item(15,1).Value = "Club Nacional ✓ - Miramar Misiones: 3:1"
$ws.Cells.Item(15,2).Value = "Club Nacional"
$ws.Cells.Item(15,3).Value = 59
$ws.Cells.Item(15,4).Value = 86
$ws.Cells.Item(15,5).Value = $null
$ws.Cells.Item(15,6).Value = 1.67
$ws.Cells.Item(15,7).Value = "✓"

# Row 16: Olympique Marseille
$ws.Cells.Item(16,1).Value = "Olympique Marseille ✓ - Le Havre AC: 6:2"
$ws.Cells.Item(16,2).Value = "Olympique Marseille"
$ws.Cells.Item(16,3).Value = 58
$ws.Cells.Item(16,4).Value = 97
$ws.Cells.Item(16,5).Value = 94
$ws.Cells.Item(16,6).Value = 1.38
$ws.Cells.Item(16,7).Value = "✓"

# Row 17: SK Brann
$ws.Cells.Item(17,1).Value = "SK Brann ✓ - FK Haugesund: 4:1"
$ws.Cells.Item(17,2).Value = "SK Brann"
$ws.Cells.Item(17,3).Value = 57
$ws.Cells.Item(17,4).Value = 100
$ws.Cells.Item(17,5).Value = $null
$ws.Cells.Item(17,6).Value = 1.11
$ws.Cells.Item(17,7).Value = "✓"

# Row 18: Club Brugge KV
$ws.Cells.Item(18,1).Value = "Oud-Heverlee Leuven - Club Brugge KV ✓: 0:1"
$ws.Cells.Item(18,2).Value = "Club Brugge KV"
$ws.Cells.Item(18,3).Value = 55
$ws.Cells.Item(18,4).Value = 86
$ws.Cells.Item(18,5).Value = $null
$ws.Cells.Item(18,6).Value = 1.67
$ws.Cells.Item(18,7).Value = "✓"
